$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) First three rows: replace their single value with "0M"
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# 2) Insert ten new rows (each holding one value) right after the third row,
#    i.e. before what is currently row 4. Insert in reverse order, always
#    anchored on the same "row 4" object, so the final order matches the
#    desired sequence.
$newRowValues = @(
    "1752",
    "0.00002",
    "0.00014",
    "0.00004",
    "0.00001",
    "0.00004",
    "0.00004",
    "0.00005",
    "0.07448",
    "100.0"
)
$anchorRow = $t.Rows.Item(4)
for ($i = $newRowValues.Count - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($anchorRow)
    $newRow.Cells.Item(1).Range.Text = $newRowValues[$i]
}

# 3) The three formerly tab-delimited "summary" rows near the end of the
#    table (each packed 10 values into one run via w:tab) collapse down to
#    a single value, re-using the values that used to live in rows 1-3.
$rowCount = $t.Rows.Count
$t.Rows.Item($rowCount - 2).Cells.Item(1).Range.Text = "99.98"
$t.Rows.Item($rowCount - 1).Cells.Item(1).Range.Text = "0.07"
$t.Rows.Item($rowCount).Cells.Item(1).Range.Text = "476"
